# Initial Chat code Commit
# Adds a new "Chat_With_Carrier_TC001" automation test row to the
# "Automation Tests" sheet, corrects a few Status values from Yes/YES to
# No, widens column A, and grows the sheet's used range / AutoFilter /
# _FilterDatabase defined name from A1:F27 to A1:F28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Automation Tests")

# --- Correct existing Status cells (col C) from Yes/YES to No ---------
$ws.Range("C20").Value = "No"
$ws.Range("C25").Value = "No"
$ws.Range("C26").Value = "No"

# Row 25's Status cell picks up the same vertical-center alignment used
# by the surrounding rows once edited.
$ws.Range("C25").VerticalAlignment = -4108

# --- Append the new test-case row (row 28) ------------------------------
# Set values in the same order the new unique strings were introduced
# (Chat Validated successfully, Validate Chat, Chat_With_Carrier_TC001)
# so they land in the shared-string table in that order.
$ws.Range("D28").Value = "Chat Validated successfully"
$ws.Range("B28").Value = "Validate Chat"
$ws.Range("A28").Value = "Chat_With_Carrier_TC001"
$ws.Range("C28").Value = "Yes"

# Match styling used by the neighbouring rows for the new cells.
$ws.Range("B28").WrapText = $true
$ws.Range("C28").VerticalAlignment = -4108
$ws.Range("D28").VerticalAlignment = -4108

# --- Widen column A to fit the new, longer test-case name --------------
$ws.Columns.Item(1).ColumnWidth = 32

# --- Grow AutoFilter to cover the new row -------------------------------
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("A1:F28").AutoFilter() | Out-Null

# AutoFilter() alone doesn't repoint the hidden _FilterDatabase defined
# name, so update it explicitly to track the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='" + $ws.Name + "'!`$A`$1:`$F`$28"
    }
}

# --- Update sheet view so row 28 is visible/selected --------------------
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B27").Select() | Out-Null
